$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = 0.4734133790737564
$ws.Range("D3").Value = 0.8267581475128645
$ws.Range("E3").Value = 0.9845626072041166
$ws.Range("H3").Value = 0.5669345898004434
$ws.Range("I3").Value = 0.1152489023255261
$ws.Range("J3").Value = 0.3739279588336192
$ws.Range("K3").Value = 1056.445969125214

$ws.Range("Q3").Value = 39
$ws.Range("R3").Value = 138
$ws.Range("S3").Value = 511
$ws.Range("T3").Value = 1114
$ws.Range("U3").Value = 1726
$ws.Range("V3").Value = 6594
$ws.Range("W3").Value = 6495
$ws.Range("X3").Value = 6122
$ws.Range("Y3").Value = 5519
$ws.Range("Z3").Value = 4907

$ws.Range("AF3").Value = 0.99412
$ws.Range("AG3").Value = 0.979195
$ws.Range("AH3").Value = 0.922961
$ws.Range("AI3").Value = 0.832052
$ws.Range("AJ3").Value = 0.7397860000000001
